$wb = $excel.ActiveWorkbook

# Add the new worksheet "Results0" right after "Sheet1"
$sheet1 = $wb.Worksheets.Item("Sheet1")
$ws = $wb.Worksheets.Add($null, $sheet1)
$ws.Name = "Results0"

# Header row (row 1) - column headers as strings
$headers = @(
    "input_Latitude",
    "input_Longitude",
    "input_PV_location",
    "input_PV_installed_cost",
    "input_Wind_installed_cost",
    "input_Site_electric_load",
    "input_Site_building_type",
    "input_Site_roofspace",
    "input_Site_landspace",
    "input_Site_NEM_limit",
    "input_Site_net_billing_rate",
    "input_Site_electricity_cost_per_kwh",
    "input_Site_demand_charge_cost_per_kw",
    "output_PV_size",
    "output_PV_energy_lcoe",
    "output_PV_energy_exported",
    "output_PV_energy_curtailed",
    "output_Wind_size",
    "output_Wind_energy_lcoe",
    "output_Wind_energy_exported",
    "output_Wind_energy_curtailed",
    "output_Grid_Electricity_Supplied_kWh_annual",
    "output_npv",
    "output_lcc"
)

# Data row (row 2) - values
$values = @(
    0,
    0,
    0,
    0.0,
    0.0,
    0.0,
    0,
    0.0,
    0.0,
    0.0,
    0.0,
    0.0,
    0.0,
    455.0,
    0.0,
    275124.0,
    0.0,
    61.0,
    0.0,
    24737.0,
    0.0,
    20155.0,
    2044361.31,
    1520728.71
)

for ($i = 0; $i -lt $headers.Length; $i++) {
    $col = $i + 1
    $ws.Cells.Item(1, $col).Value = $headers[$i]
    $ws.Cells.Item(2, $col).Value = $values[$i]
}

# Restore Sheet1 as the active/selected tab (matches target workbook state)
$sheet1.Activate()

Write-Output "edit applied"
